$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.057.09"
$ws.Range("E2").Value = "  -3.60%  "

$ws.Range("D3").Value = "3.518.65"
$ws.Range("E3").Value = "  -4.24%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.11"
$ws.Range("E5").Value = "  -5.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.15"
$ws.Range("E6").Value = "  -4.78%  "

$ws.Range("D7").Value = "3.517.83"
$ws.Range("E7").Value = "  -4.21%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -4.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.88"
$ws.Range("E11").Value = "  -4.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -3.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -5.36%  "

$ws.Range("D14").Value = "4.114.04"
$ws.Range("E14").Value = "  -4.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.54"
$ws.Range("E15").Value = "  -3.43%  "

$ws.Range("D16").Value = "3.521.91"
$ws.Range("E16").Value = "  -4.74%  "

$ws.Range("D17").Value = "66.932.33"
$ws.Range("E17").Value = "  -3.79%  "

$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("E19").Value = "  -3.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.33"
$ws.Range("E20").Value = "  -4.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.29"
$ws.Range("E21").Value = "  -5.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.12"
$ws.Range("E22").Value = "  -9.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.629"
$ws.Range("E23").Value = "  -2.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.76"
$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("D26").Value = "3.660.62"
$ws.Range("E26").Value = "  -4.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("E27").Value = "  -2.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  -6.55%  "

$ws.Range("E29").Value = "  -10.07%  "

$ws.Range("E30").Value = "  -3.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -3.35%  "

$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.67"
$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.159"
$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.14"
$ws.Range("E35").Value = "  -4.67%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  -7.16%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.512.88"
$ws.Range("E37").Value = "  -4.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.02"
$ws.Range("E38").Value = "  -5.11%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.29"
$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.14"
$ws.Range("E42").Value = "  -3.06%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.56"
$ws.Range("E43").Value = "  -5.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0859"
$ws.Range("E44").Value = "  -4.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.889"
$ws.Range("E45").Value = "  -3.91%  "

$ws.Range("E46").Value = "  -4.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.01"
$ws.Range("E47").Value = "  -6.96%  "

$ws.Range("E48").Value = "  -5.64%  "

$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.55"
$ws.Range("E50").Value = "  -3.42%  "

$ws.Range("E51").Value = "  -4.73%  "
